# Import file excel from FE
# Clear the cells that were removed in the FE import (duplicate/obsolete rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("A7").ClearContents()

$ws.Range("G5").Select()
